{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\n\nconst words = [\"Week\", \" \", \"7\", \" \", \"Reading\", \" \", \"Guide\", \" \", \"Part\", \" \", \"2:\", \" \", \"Hypothesis\", \" \", \"Tests\"];\n\n// Clear paragraph content first\nconst r = titlePara.getRange();\nr.insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\nfor (const w of words) {\n  const ins = titlePara.insertText(w, Word.InsertLocation.end);\n  ins.font.bold = true;\n  ins.font.bold = false;\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\nWrite-Output $d.Paragraphs.Count\n"}
